# LOM3058.xlsx edit: fill in the "Objetivos", "Programa" and "Bibliografia"
# sections with real content, and add a dedicated "Docentes responsáveis:"
# row holding the professor's name (previously parked under "Objetivos:").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the old "Programa resumido:" row (row 13) so the
# professor name can move out from under "Objetivos:" into its own row
# below "Docentes responsáveis:". Everything from row 13 down shifts to 14+.
$ws.Rows.Item(13).EntireRow.Insert()

# Row 10 "Objetivos:" - was holding the professor name; now holds the real
# Portuguese objectives text.
$ws.Range("B10").Value = "Apresentar ao estudante conceitos gerais sobre a síntese de polímeros, destacando as principais vias usadas para a geração de materiais poliméricos na atualidade. Capacitar o estudante para relacionar a síntese com a estrutura, o comportamento e a utilização de polímeros."
$ws.Range("C10").Value = "Apresentar ao estudante conceitos gerais sobre a síntese de polímeros, destacando as principais vias usadas para a geração de materiais poliméricos na atualidade. Capacitar o estudante para relacionar a síntese com a estrutura, o comportamento e a utilização de polímeros."

# New row 13 (blank after insert), under "Docentes responsáveis:" - holds
# the professor name that used to live under "Objetivos:".
$ws.Range("B13").Value = "1033242 - Fábio Herbst Florenzano"
$ws.Range("C13").Value = "1033242 - Fábio Herbst Florenzano"

# Row 14 "Programa resumido:" (was row 13) - was a placeholder "Semestral",
# now the real short syllabus text in Portuguese.
$ws.Range("B14").Value = "Fundamentos: massa molar média e conceitos gerais sobre química de polímeros. Polímeros de condensação e de adição. Polimerização em cadeia. Polimerização em etapas. Métodos de Polimerização. Modificação e degradação de polímeros e seu retardo."
$ws.Range("C14").Value = "Fundamentos: massa molar média e conceitos gerais sobre química de polímeros. Polímeros de condensação e de adição. Polimerização em cadeia. Polimerização em etapas. Métodos de Polimerização. Modificação e degradação de polímeros e seu retardo."

# Row 16 "Programa:" (was row 15) - was a placeholder date, now the real
# long syllabus text in Portuguese.
$ws.Range("B16").Value = "Fundamentos: massa molar média e conceitos gerais sobre química de polímeros. Polímeros de condensação e de adição: conceitos gerais. Polimerização em Etapas. Polimerização em cadeia: via radical, aniônica e catiônica. Polimerização via desativação reversível de radicais. Polimerização estéreo-específica: Ziegler-Natta e outras. Métodos de polimerização: batelada, solução, suspensão, emulsão e interfacial. Modificação de polímeros: reticulação e vulcanização; modificações em aromáticos, hidrólise e outras. Química da degradação de polímeros: processos gerais e métodos de controle."
$ws.Range("C16").Value = "Fundamentos: massa molar média e conceitos gerais sobre química de polímeros. Polímeros de condensação e de adição: conceitos gerais. Polimerização em Etapas. Polimerização em cadeia: via radical, aniônica e catiônica. Polimerização via desativação reversível de radicais. Polimerização estéreo-específica: Ziegler-Natta e outras. Métodos de polimerização: batelada, solução, suspensão, emulsão e interfacial. Modificação de polímeros: reticulação e vulcanização; modificações em aromáticos, hidrólise e outras. Química da degradação de polímeros: processos gerais e métodos de controle."

# Row 19 "Método:" (was row 18) - was holding the professor name by
# mistake; now holds the real evaluation method text.
$ws.Range("B19").Value = "Provas escritas envolvendo o conteúdo teórico ministrado em sala de aula."
$ws.Range("C19").Value = "Provas escritas envolvendo o conteúdo teórico ministrado em sala de aula."

# Row 20 "Critério:" (was row 19) - now holds the pass/fail criteria text
# that used to (incorrectly) sit under "Norma de recuperação:".
$ws.Range("B20").Value = "Duas avaliações, sendo que a nota final corresponde à média aritmética das duas provas. Os alunos que apresentarem média igual ou superior a 5 estarão aprovados, enquanto que aqueles que tiverem média inferior a 3 estarão reprovados. Alunos com notas finais situadas no intervalo de 3 a 5 serão encaminhados à recuperação."
$ws.Range("C20").Value = "Duas avaliações, sendo que a nota final corresponde à média aritmética das duas provas. Os alunos que apresentarem média igual ou superior a 5 estarão aprovados, enquanto que aqueles que tiverem média inferior a 3 estarão reprovados. Alunos com notas finais situadas no intervalo de 3 a 5 serão encaminhados à recuperação."

# Row 21 "Norma de recuperação:" (was row 20) - now holds the actual
# recovery-exam rules, previously mislabeled under "Bibliografia:".
$ws.Range("B21").Value = "O aluno será submetido a um programa de estudos destinado a rever o conteúdo abordado na disciplina. Ao final deste período será aplicada uma nova avaliação. A nota final do aluno será a média aritmética desta avaliação com a nota anteriormente obtida, estando aprovados os alunos que tiverem nota final igual ou superior a 5."
$ws.Range("C21").Value = "O aluno será submetido a um programa de estudos destinado a rever o conteúdo abordado na disciplina. Ao final deste período será aplicada uma nova avaliação. A nota final do aluno será a média aritmética desta avaliação com a nota anteriormente obtida, estando aprovados os alunos que tiverem nota final igual ou superior a 5."

# Row 22 "Bibliografia:" (was row 21) - now holds the real bibliography
# list.
$ws.Range("B22").Value = "G. ODIAN Principles of Polymerization, 3rd Edition, New York: Wiley-Interscience, 1991.`nF. W. Billmeyer. Textbook of Polymer Chemistry, 3rd edition, New York: Wiley-Interscience, 1984.`nC. E. Carraher. Introduction to Polymer Chemistry, 1st Edition, Boca Raton: Taylor and Francis, 2010.`nS. V. Canevarolo. Ciência dos Polímeros: um texto básico para Engenheiros e Tecnólogos, 2ª. edição, São Paulo: Artliber, 2006."
$ws.Range("C22").Value = "G. ODIAN Principles of Polymerization, 3rd Edition, New York: Wiley-Interscience, 1991.`nF. W. Billmeyer. Textbook of Polymer Chemistry, 3rd edition, New York: Wiley-Interscience, 1984.`nC. E. Carraher. Introduction to Polymer Chemistry, 1st Edition, Boca Raton: Taylor and Francis, 2010.`nS. V. Canevarolo. Ciência dos Polímeros: um texto básico para Engenheiros e Tecnólogos, 2ª. edição, São Paulo: Artliber, 2006."
